# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (only wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     colours (wired to the Slide Master / deck)
#
# The authored edit swaps the two themes' content, so the Slide Master (and
# therefore every slide) ends up using the stock "Office Theme" palette
# instead of "Integral". The 12-colour theme scheme is the only thing that
# differs between the two theme parts (fonts / fill-format scheme are
# byte-identical already), so re-pointing the Slide Master's theme colours
# at the Office Theme palette reproduces the meaningful, visible effect of
# the swap through the supported PowerPoint object model.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$scheme = $m.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" colours that used to live in
# theme1.xml, expressed as COM RGB() (0x00BBGGRR) values.
$targetColors = @(
    @{ Index = 1;  Hex = '000000'; RGB = 0        },  # dk1
    @{ Index = 2;  Hex = 'FFFFFF'; RGB = 16777215  },  # lt1
    @{ Index = 3;  Hex = '44546A'; RGB = 6968388   },  # dk2
    @{ Index = 4;  Hex = 'E7E6E6'; RGB = 15132391  },  # lt2
    @{ Index = 5;  Hex = '5B9BD5'; RGB = 13998939  },  # accent1
    @{ Index = 6;  Hex = 'ED7D31'; RGB = 3243501   },  # accent2
    @{ Index = 7;  Hex = 'A5A5A5'; RGB = 10855845  },  # accent3
    @{ Index = 8;  Hex = 'FFC000'; RGB = 49407     },  # accent4
    @{ Index = 9;  Hex = '4472C4'; RGB = 12874308  },  # accent5
    @{ Index = 10; Hex = '70AD47'; RGB = 4697456   },  # accent6
    @{ Index = 11; Hex = '0563C1'; RGB = 12673797  },  # hlink
    @{ Index = 12; Hex = '954F72'; RGB = 7491477   }   # folHlink
)

foreach ($entry in $targetColors) {
    $scheme.Item($entry.Index).RGB = $entry.RGB
}
